# Netflix.xlsx edit: extend the cost/participant ratio + year roll-forward
# formulas from column BU out to column CG (12 more months), break the
# BU column out of each shared-formula block (it references its own column
# rather than the +12 offset like the rest of the block), mark Rida's
# month as paid (checkmark) for gennaio-marzo, and update Rida's running
# balance plus the last-used-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Row 2 ("✔️ / ratio" row)
# ---------------------------------------------------------------------
# BU2 no longer follows the shared +12 offset pattern used by Z2:BT2 -
# it now references its own column, so it must become a standalone formula.
$ws.Range("BU2").Formula = "=ROUND(BU8/BU9,1)"
# New columns BV2:CG2 continue the ratio pattern as a fresh shared block.
$ws.Range("BV2:CG2").Formula = "=ROUND(BV8/BV9,1)"

# ---------------------------------------------------------------------
# Row 3 (Rida's row)
# ---------------------------------------------------------------------
# Rida's payment is now marked "paid" (checkmark) for gennaio-marzo 2026
# (AI3:AK3), so those cells stop being part of the formula block.
$ws.Range("AI3:AK3").Value = "✔️"
# The shared block's new anchor starts at AL3 (same formula text/offset).
$ws.Range("AL3:BM3").Formula = "=ROUND(AX8/AX9,1)"
$ws.Range("BN3:BU3").Formula = "=ROUND(BZ8/BZ9,1)"
# BU3 breaks the pattern (self-referencing column) just like BU2/BU4.
$ws.Range("BU3").Formula = "=ROUND(BU8/BU9,1)"
# New columns BV3:CG3.
$ws.Range("BV3:CG3").Formula = "=ROUND(BV8/BV9,1)"

# ---------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------
$ws.Range("AD4:AT4").Formula = "=ROUND(AP8/AP9,1)"
$ws.Range("AU4:BU4").Formula = "=ROUND(BG8/BG9,1)"
# BU4 breaks the pattern (self-referencing column) just like BU2/BU3.
$ws.Range("BU4").Formula = "=ROUND(BU8/BU9,1)"
# New columns BV4:CG4.
$ws.Range("BV4:CG4").Formula = "=ROUND(BV8/BV9,1)"

# ---------------------------------------------------------------------
# Row 5 (year roll-forward)
# ---------------------------------------------------------------------
$ws.Range("C5:AH5").Formula = "=IF(N1=12,B5+1,B5)"
$ws.Range("AI5:BN5").Formula = "=IF(AT1=12,AH5+1,AH5)"
$ws.Range("BO5:BU5").Formula = "=IF(BZ1=12,BN5+1,BN5)"
# New columns BV5:CG5 extend the year roll-forward one cell at a time.
$ws.Range("BV5").Formula = "=IF(CG1=12,BU5+1,BU5)"
$ws.Range("BW5").Formula = "=IF(CH1=12,BV5+1,BV5)"
$ws.Range("BX5").Formula = "=IF(CI1=12,BW5+1,BW5)"
$ws.Range("BY5").Formula = "=IF(CJ1=12,BX5+1,BX5)"
$ws.Range("BZ5").Formula = "=IF(CK1=12,BY5+1,BY5)"
$ws.Range("CA5").Formula = "=IF(CL1=12,BZ5+1,BZ5)"
$ws.Range("CB5").Formula = "=IF(CM1=12,CA5+1,CA5)"
$ws.Range("CC5").Formula = "=IF(CN1=12,CB5+1,CB5)"
$ws.Range("CD5").Formula = "=IF(CO1=12,CC5+1,CC5)"
$ws.Range("CE5").Formula = "=IF(CP1=12,CD5+1,CD5)"
$ws.Range("CF5").Formula = "=IF(CQ1=12,CE5+1,CE5)"
$ws.Range("CG5").Formula = "=IF(CR1=12,CF5+1,CF5)"

# ---------------------------------------------------------------------
# Rida paid 200 EUR on 25/12/2025: her running balance (B17) grows from
# 0.9 to 185.9 (the remainder is kept aside for future price changes).
# ---------------------------------------------------------------------
$ws.Range("B17").Value = 185.9

# Leave the selection where the user's last edit was (B18).
$ws.Range("B18").Select() | Out-Null
